# Updates the "cryptos" price table (Sheet1, rows 2-51) with refreshed
# Price / Volume(1h) figures, and two coin re-rankings:
#   - rows 16/17 swap: TRON/Chainlink -> Chainlink/TRON
#   - rows 49/50 swap: Monero/dogwifhat -> dogwifhat/Monero
#   - row 51: THORChain replaced by LidoDAOToken
#
# D (Price) and E (Volume) columns are stored as plain text in the
# workbook (e.g. "66.220.74", "  +2.01%  "), so a leading apostrophe
# forces Excel to keep them as text instead of re-parsing them as
# numbers/percentages (which would silently rewrite "1.00" -> "1",
# "0.120" -> "0.12", "2.15%" -> 0.0215, etc). Resetting the style back
# to "Normal" afterwards avoids leaving a stray Text/quote-prefix
# number format on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'66.362.23"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +2.15%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'3.417.54"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +0.82%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("D4").Value = "'0.998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  -0.24%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'566.90"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +1.36%  "
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = "'178.37"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +2.47%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").Value = "'0.632"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +1.32%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = "'3.403.01"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +0.70%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = "'0.999"
$ws.Range("D9").Style = "Normal"

$ws.Range("D10").Value = "'0.177"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +4.35%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Value = "'0.638"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +1.07%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = "'54.49"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -0.12%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = "'0.0000281"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +0.59%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'9.32"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +2.15%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'3.934.74"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +0.11%  "
$ws.Range("E15").Style = "Normal"

# Row 16: coin identity change
$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").Value = "'18.32"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +0.15%  "
$ws.Range("E16").Style = "Normal"

# Row 17: coin identity change
$ws.Range("B17").Value = "TRON"
$ws.Range("C17").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D17").Value = "'0.120"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +0.80%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = "'3.410.05"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +0.64%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Value = "'66.063.28"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +1.88%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").Value = "'11.97"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +1.23%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("E21").Value = "'  +1.18%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value = "'466.19"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -1.06%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "'4.97"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -0.09%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("E24").Value = "'  +9.90%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = "'90.08"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +3.48%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").Value = "'4.13"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -0.24%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("D27").Value = "'2.94"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +0.72%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("D28").Value = "'10.75"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -0.36%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("D29").Value = "'8.82"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +0.53%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("D30").Value = "'31.39"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +1.29%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("D31").Value = "'6.69"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +0.06%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("D32").Value = "'11.55"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +0.28%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("D33").Value = "'581.86"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +1.59%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("D34").Value = "'62.57"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +1.83%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("E35").Value = "'  +0.78%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("E36").Value = "'  -0.01%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("E37").Value = "'  +2.46%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("E38").Value = "'  +0.49%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("D39").Value = "'36.43"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +2.12%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").Value = "'0.382"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +3.18%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").Value = "'0.0₃0761"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +1.43%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("D42").Value = "'3.131.32"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +1.24%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Value = "'2.88"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +0.74%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = "'0.0421"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +1.41%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = "'2.49"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +1.00%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("E46").Value = "'  -0.44%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("D47").Value = "'3.19"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +0.41%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").Value = "'1.00"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +0.13%  "
$ws.Range("E48").Style = "Normal"

# Row 49: coin identity change
$ws.Range("B49").Value = "dogwifhat"
$ws.Range("C49").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D49").Value = "'2.61"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +11.82%  "
$ws.Range("E49").Style = "Normal"

# Row 50: coin identity change
$ws.Range("B50").Value = "Monero"
$ws.Range("C50").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D50").Value = "'141.30"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +2.35%  "
$ws.Range("E50").Style = "Normal"

# Row 51: coin identity change
$ws.Range("B51").Value = "LidoDAOToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D51").Value = "'3.18"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +10.15%  "
$ws.Range("E51").Style = "Normal"
